$d = $word.ActiveDocument

# Collapse the split "<id>...</id>" runs (three runs: "<id>", the bare id
# text, "</id>") into a single run containing the full "<id>...</id>"
# text, for each of the two ids that were re-downloaded.
$ids = @("p112v_1", "p113r_1")

foreach ($idVal in $ids) {
    $old = "<id>" + $idVal + "</id>"
    $new = "<id>" + $idVal + "</id>"
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Output "Replaced '$idVal': $found"
}
